$wb = $excel.ActiveWorkbook

# Column-letter -> 1-based index map for columns H..N used below.
$colIndex = @{ "H" = 8; "I" = 9; "J" = 10; "K" = 11; "L" = 12; "M" = 13; "N" = 14 }

# Each entry is one row whose currentAveragePrice / LevePrice / LeveProfit
# columns (H-N) were refreshed by the scheduled pricing-data runner.
# A value of $null means the cell is cleared (no profit computable for that column).
$updates = @(
    @{ Sheet = "ALC"; Row = 38; Cells = @{ "H" = 674.8333; "I" = 159.8; "J" = 3250; "K" = 479.4; "L" = 9750; "M" = -107.4; "N" = -10494 } },
    @{ Sheet = "ALC"; Row = 58; Cells = @{ "H" = 967.7; "I" = 460.25; "J" = 2997.5; "K" = 1380.75; "L" = 8992.5; "M" = -1230.75; "N" = -9292.5 } },
    @{ Sheet = "ALC"; Row = 62; Cells = @{ "H" = 4994.6; "I" = 4994.6; "J" = 0; "K" = 4994.6; "L" = 0; "M" = -4370.6 } },
    @{ Sheet = "ALC"; Row = 65; Cells = @{ "H" = 4994.6; "I" = 4994.6; "J" = 0; "K" = 24973; "L" = 0; "M" = -21853 } },
    @{ Sheet = "ALC"; Row = 69; Cells = @{ "H" = 0; "I" = 0; "J" = 0; "K" = 0; "L" = 0; "M" = $null } },
    @{ Sheet = "ALC"; Row = 72; Cells = @{ "H" = 0; "I" = 0; "J" = 0; "K" = 0; "L" = 0; "M" = $null } },
    @{ Sheet = "ALC"; Row = 86; Cells = @{ "H" = 3885.6875; "I" = 3679.3635; "J" = 4339.6; "K" = 3679.3635; "L" = 4339.6; "M" = -2556.3635; "N" = -6585.6 } },
    @{ Sheet = "ALC"; Row = 89; Cells = @{ "H" = 3885.6875; "I" = 3679.3635; "J" = 4339.6; "K" = 18396.8175; "L" = 21698; "M" = -12780.8175; "N" = -32930 } },
    @{ Sheet = "ALC"; Row = 103; Cells = @{ "H" = 435.85715; "I" = 423.5; "J" = 440.8; "K" = 1270.5; "L" = 1322.4; "M" = -684.5; "N" = -2494.4 } },
    @{ Sheet = "ALC"; Row = 138; Cells = @{ "H" = 1763.3334; "I" = 580; "J" = 2000; "K" = 1740; "L" = 6000; "M" = 3400 } },
    @{ Sheet = "ARM"; Row = 122; Cells = @{ "H" = 0; "I" = 0; "J" = 0; "K" = 0; "L" = 0; "M" = $null; "N" = $null } },
    @{ Sheet = "BSM"; Row = 75; Cells = @{ "H" = 22815.572; "I" = 5377.25; "J" = 46066.668; "K" = 5377.25; "L" = 46066.668; "M" = -4441.25; "N" = -47938.668 } },
    @{ Sheet = "BSM"; Row = 78; Cells = @{ "H" = 22815.572; "I" = 5377.25; "J" = 46066.668; "K" = 16131.75; "L" = 138200.004; "M" = -11451.75; "N" = -147560.004 } },
    @{ Sheet = "CRP"; Row = 55; Cells = @{ "H" = 9500; "I" = 9500; "J" = 0; "K" = 9500; "L" = 0; "M" = -9185 } },
    @{ Sheet = "CRP"; Row = 59; Cells = @{ "H" = 31883.334; "I" = 28983.334; "J" = 33333.332; "K" = 28983.334; "L" = 33333.332; "M" = -27838.334; "N" = -35623.332 } },
    @{ Sheet = "CRP"; Row = 60; Cells = @{ "H" = 21452.727; "I" = 12000; "J" = 24997.5; "K" = 12000; "L" = 24997.5; "M" = -11489; "N" = -26019.5 } },
    @{ Sheet = "CRP"; Row = 68; Cells = @{ "H" = 39998.637; "I" = 0; "J" = 39998.637; "K" = 0; "L" = 39998.637; "N" = -41496.637 } },
    @{ Sheet = "CRP"; Row = 71; Cells = @{ "H" = 39998.637; "I" = 0; "J" = 39998.637; "K" = 0; "L" = 119995.911; "N" = -127483.911 } },
    @{ Sheet = "CRP"; Row = 74; Cells = @{ "H" = 39499.145; "I" = 36494; "J" = 40000; "K" = 36494; "L" = 40000; "M" = -35620; "N" = -41748 } },
    @{ Sheet = "CRP"; Row = 77; Cells = @{ "H" = 39499.145; "I" = 36494; "J" = 40000; "K" = 109482; "L" = 120000; "M" = -105114; "N" = -128736 } },
    @{ Sheet = "CRP"; Row = 107; Cells = @{ "H" = 866.9091; "I" = 921.5; "J" = 721.3333; "K" = 921.5; "L" = 721.3333; "M" = 998.5; "N" = -4561.3333 } },
    @{ Sheet = "CRP"; Row = 132; Cells = @{ "H" = 3499.6667; "I" = 3499.5; "J" = 3500; "K" = 10498.5; "L" = 10500; "M" = -7968.5 } },
    @{ Sheet = "CRP"; Row = 134; Cells = @{ "H" = 8337.333000000001; "I" = 8337.333000000001; "J" = 0; "K" = 25011.999; "L" = 0; "M" = -22476.999 } },
    @{ Sheet = "CUL"; Row = 18; Cells = @{ "H" = 1038.3334; "I" = 1038.3334; "J" = 0; "K" = 3115.0002; "L" = 0; "M" = -2946.0002 } },
    @{ Sheet = "CUL"; Row = 47; Cells = @{ "H" = 93.666664; "I" = 93.666664; "J" = 0; "K" = 280.999992; "L" = 0; "M" = 150.000008; "N" = $null } },
    @{ Sheet = "CUL"; Row = 102; Cells = @{ "H" = 500; "I" = 500; "J" = 0; "K" = 1500; "L" = 0; "M" = 934 } },
    @{ Sheet = "CUL"; Row = 103; Cells = @{ "H" = 269.6; "I" = 309.85715; "J" = 175.66667; "K" = 929.5714499999999; "L" = 527.00001; "M" = -50.57144999999991; "N" = -2285.00001 } },
    @{ Sheet = "CUL"; Row = 120; Cells = @{ "H" = 8300; "I" = 5450; "J" = 14000; "K" = 16350; "L" = 42000; "M" = -11512 } },
    @{ Sheet = "CUL"; Row = 131; Cells = @{ "H" = 1598.5; "I" = 667.3333; "J" = 1997.5714; "K" = 2001.9999; "L" = 5992.7142; "M" = 3038.0001 } },
    @{ Sheet = "GSM"; Row = 44; Cells = @{ "H" = 25000.5; "I" = 0; "J" = 25000.5; "K" = 0; "L" = 25000.5; "N" = -26192.5 } },
    @{ Sheet = "LTW"; Row = 16; Cells = @{ "H" = 1477.8; "I" = 1477.8; "J" = 0; "K" = 1477.8; "L" = 0; "M" = -1307.8 } },
    @{ Sheet = "LTW"; Row = 22; Cells = @{ "H" = 1856.875; "I" = 1858.7646; "J" = 1852.2858; "K" = 1858.7646; "L" = 1852.2858; "M" = -1563.7646; "N" = -2442.2858 } },
    @{ Sheet = "LTW"; Row = 27; Cells = @{ "H" = 1856.875; "I" = 1858.7646; "J" = 1852.2858; "K" = 1858.7646; "L" = 1852.2858; "M" = -1751.7646; "N" = -2066.2858 } },
    @{ Sheet = "LTW"; Row = 40; Cells = @{ "H" = 6775.25; "I" = 3550.5; "J" = 10000; "K" = 3550.5; "L" = 10000; "M" = -3414.5 } },
    @{ Sheet = "LTW"; Row = 43; Cells = @{ "H" = 11000; "I" = 0; "J" = 11000; "K" = 0; "L" = 11000; "M" = $null; "N" = -11386 } },
    @{ Sheet = "LTW"; Row = 46; Cells = @{ "H" = 2422.5715; "I" = 3000; "J" = 2191.6; "K" = 3000; "L" = 2191.6; "M" = -2812; "N" = -2567.6 } },
    @{ Sheet = "LTW"; Row = 69; Cells = @{ "H" = 70163; "I" = 0; "J" = 70163; "K" = 0; "L" = 70163; "N" = -71785 } },
    @{ Sheet = "LTW"; Row = 72; Cells = @{ "H" = 70163; "I" = 0; "J" = 70163; "K" = 0; "L" = 210489; "N" = -218601 } },
    @{ Sheet = "LTW"; Row = 122; Cells = @{ "H" = 3229.8; "I" = 3229.8; "J" = 0; "K" = 9689.400000000001; "L" = 0; "M" = -7239.400000000001 } },
    @{ Sheet = "LTW"; Row = 132; Cells = @{ "H" = 10527.308; "I" = 9050.714; "J" = 12250; "K" = 27152.142; "L" = 36750; "M" = -24622.142; "N" = -41810 } },
    @{ Sheet = "LTW"; Row = 136; Cells = @{ "H" = 2987.5386; "I" = 3012.0908; "J" = 2852.5; "K" = 9036.2724; "L" = 8557.5; "M" = -6486.2724 } },
    @{ Sheet = "WVR"; Row = 62; Cells = @{ "H" = 4779.5; "I" = 4759.909; "J" = 4995; "K" = 4759.909; "L" = 4995; "M" = -4135.909; "N" = -6243 } },
    @{ Sheet = "WVR"; Row = 65; Cells = @{ "H" = 4779.5; "I" = 4759.909; "J" = 4995; "K" = 23799.545; "L" = 24975; "M" = -20679.545; "N" = -31215 } },
    @{ Sheet = "WVR"; Row = 81; Cells = @{ "H" = 398; "I" = 398; "J" = 0; "K" = 796; "L" = 0; "M" = 265 } },
    @{ Sheet = "WVR"; Row = 84; Cells = @{ "H" = 398; "I" = 398; "J" = 0; "K" = 3980; "L" = 0; "M" = 1324 } },
    @{ Sheet = "WVR"; Row = 107; Cells = @{ "H" = 254.5; "I" = 262.2857; "J" = 200; "K" = 786.8571000000001; "L" = 600; "M" = 1133.1429 } },
    @{ Sheet = "WVR"; Row = 136; Cells = @{ "H" = 1321.5714; "I" = 1321.5714; "J" = 0; "K" = 3964.7142; "L" = 0; "M" = -1414.7142 } }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    foreach ($col in $u.Cells.Keys) {
        $colNum = $colIndex[$col]
        $ws.Cells.Item($u.Row, $colNum).Value = $u.Cells[$col]
    }
}
